$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "29.201.59"
$ws.Range("E2").Value = "  -0.91%  "

# Row 3
Set-TextValue "D3" "1.866.82"
$ws.Range("E3").Value = "  -0.48%  "

# Row 4
Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
Set-TextValue "D5" "0.7097"
$ws.Range("E5").Value = "  -0.78%  "

# Row 6
$ws.Range("E6").Value = "  +0.12%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
Set-TextValue "D8" "0.3110"
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
Set-TextValue "D9" "0.07648"
$ws.Range("E9").Value = "  -3.75%  "

# Row 10
Set-TextValue "D10" "24.69"
$ws.Range("E10").Value = "  -2.59%  "

# Row 11
Set-TextValue "D11" "0.08371"
$ws.Range("E11").Value = "  +1.13%  "

# Row 12
Set-TextValue "D12" "1.864.49"
$ws.Range("E12").Value = "  -0.70%  "

# Row 13
Set-TextValue "D13" "5.226"
$ws.Range("E13").Value = "  -0.99%  "

# Row 14
Set-TextValue "D14" "0.7104"
$ws.Range("E14").Value = "  -2.55%  "

# Row 15
$ws.Range("E15").Value = "  +0.18%  "

# Row 16
Set-TextValue "D16" "29.211.22"
$ws.Range("E16").Value = "  -0.91%  "

# Row 17
$ws.Range("E17").Value = "  +0.41%  "

# Row 18
Set-TextValue "D18" "243.47"
$ws.Range("E18").Value = "  -0.80%  "

# Row 19
Set-TextValue "D19" "0.000007831"
$ws.Range("E19").Value = "  -0.60%  "

# Row 20
Set-TextValue "D20" "2.115.56"
$ws.Range("E20").Value = "  +0.38%  "

# Row 21
$ws.Range("E21").Value = "  -1.96%  "

# Row 22
Set-TextValue "D22" "0.9995"
$ws.Range("E22").Value = "  +0.00%  "

# Row 23
Set-TextValue "D23" "7.856"
$ws.Range("E23").Value = "  -1.46%  "

# Row 24
$ws.Range("E24").Value = "  -0.01%  "

# Row 25
Set-TextValue "D25" "0.1629"
$ws.Range("E25").Value = "  +1.49%  "

# Row 26
Set-TextValue "D26" "163.32"
$ws.Range("E26").Value = "  -0.08%  "

# Row 27
$ws.Range("E27").Value = "  -1.14%  "

# Row 28
Set-TextValue "D28" "18.51"
$ws.Range("E28").Value = "  +1.06%  "

# Row 29
Set-TextValue "D29" "1.507"
$ws.Range("E29").Value = "  +1.05%  "

# Row 30
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D30" "4.397"
$ws.Range("E30").Value = "  +0.26%  "

# Row 31
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D31" "1.307"
$ws.Range("E31").Value = "  -3.80%  "

# Row 32
Set-TextValue "D32" "4.245"
$ws.Range("E32").Value = "  +3.16%  "

# Row 33
Set-TextValue "D33" "0.05141"
$ws.Range("E33").Value = "  -2.33%  "

# Row 34
Set-TextValue "D34" "0.7982"
$ws.Range("E34").Value = "  +9.75%  "

# Row 35
Set-TextValue "D35" "1.911"
$ws.Range("E35").Value = "  -2.27%  "

# Row 36
$ws.Range("E36").Value = "  -2.65%  "

# Row 37
Set-TextValue "D37" "2.687"
$ws.Range("E37").Value = "  +0.31%  "

# Row 38
$ws.Range("E38").Value = "  -0.73%  "

# Row 39
Set-TextValue "D39" "2.706"
$ws.Range("E39").Value = "  +0.11%  "

# Row 40
Set-TextValue "D40" "1.156.53"
$ws.Range("E40").Value = "  -5.40%  "

# Row 41
Set-TextValue "D41" "6.327"
$ws.Range("E41").Value = "  +3.33%  "

# Row 42
Set-TextValue "D42" "0.8973"
$ws.Range("E42").Value = "  -1.52%  "

# Row 43
Set-TextValue "D43" "73.26"
$ws.Range("E43").Value = "  -0.78%  "

# Row 44
Set-TextValue "D44" "0.9995"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
Set-TextValue "D45" "103.31"
$ws.Range("E45").Value = "  +1.12%  "

# Row 46
Set-TextValue "D46" "2.012.66"
$ws.Range("E46").Value = "  +0.02%  "

# Row 47
Set-TextValue "D47" "0.5177"
$ws.Range("E47").Value = "  -2.07%  "

# Row 48
Set-TextValue "D48" "1.778"
$ws.Range("E48").Value = "  -1.10%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "9.333"
$ws.Range("E49").Value = "  -0.14%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D50" "0.00000000119"
$ws.Range("E50").Value = "  -0.98%  "

# Row 51
Set-TextValue "D51" "0.4291"
$ws.Range("E51").Value = "  -0.67%  "
